# Update the quarterly report dates on the "Reporte de Formatos" sheet
# (row 8) from Q3 2021 to Q4 2021, and move the saved selection to C10
# (matching the author's last cursor position when they saved the file).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Periodo que se informa: 2021-07-01 / 2021-09-30  ->  2021-10-01 / 2021-12-31
$ws.Range("B8").Value = 44470
$ws.Range("C8").Value = 44561

# Fecha de validación / Fecha de actualización: 2021-10-11 -> 2022-01-10
$ws.Range("N8").Value = 44571
$ws.Range("O8").Value = 44571

# Leave the cursor on C10, matching the saved view in the updated file.
$ws.Range("C10").Select()
